# Updates the "cryptos" price/volume table to the latest scrape snapshot.
# Mirrors the GitHub Actions commit "Updated cryptos list ... with GitHub Actions":
#  - refreshed Price (col D) / Volume(1h) (col E) values for most rows
#  - re-ranked three coin pairs/triples, moving each row's Coin/Link/Price/Volume
#    together (rows 27-28 swap Dai/RenderToken; rows 41-43 rotate Arweave/
#    dogwifhat/Kaspa) to reflect their new relative ranking
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    # Force text storage so numeric-looking strings (e.g. "1.00", "0.999")
    # are not coerced to numbers, then drop back to the default/unstyled
    # format so no stray cell style is introduced.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2: 'Bitcoin'
Set-TextCell 'D2' '65.857.92'
Set-TextCell 'E2' '  +1.38%  '

# Row 3: 'Ethereum'
Set-TextCell 'D3' '2.954.97'
Set-TextCell 'E3' '  -1.16%  '

# Row 4: 'TetherUSD'
Set-TextCell 'E4' '  +0.07%  '

# Row 5: 'BNB'
Set-TextCell 'D5' '567.34'
Set-TextCell 'E5' '  -2.46%  '

# Row 6: 'Solana'
Set-TextCell 'D6' '160.53'
Set-TextCell 'E6' '  +4.63%  '

# Row 7: 'USDC'
Set-TextCell 'D7' '0.999'
Set-TextCell 'E7' '  -0.13%  '

# Row 8: 'XRP'
Set-TextCell 'E8' '  +1.55%  '

# Row 9: 'LidoStakedEther'
Set-TextCell 'D9' '2.951.40'
Set-TextCell 'E9' '  -1.18%  '

# Row 10: 'Toncoin'
Set-TextCell 'D10' '6.73'
Set-TextCell 'E10' '  -3.35%  '

# Row 11: 'Dogecoin'
Set-TextCell 'E11' '  -1.64%  '

# Row 12: 'Cardano'
Set-TextCell 'D12' '0.454'
Set-TextCell 'E12' '  +1.54%  '

# Row 14: 'Avalanche'
Set-TextCell 'D14' '34.21'
Set-TextCell 'E14' '  -0.24%  '

# Row 15: 'TRON'
Set-TextCell 'E15' '  -0.71%  '

# Row 16: 'WrappedBTC'
Set-TextCell 'D16' '65.893.30'
Set-TextCell 'E16' '  +1.49%  '

# Row 17: 'WrappedliquidstakedEther2.0'
Set-TextCell 'D17' '3.446.18'
Set-TextCell 'E17' '  -0.96%  '

# Row 18: 'Polkadot'
Set-TextCell 'D18' '6.93'
Set-TextCell 'E18' '  +0.21%  '

# Row 19: 'WrappedEther'
Set-TextCell 'D19' '2.953.53'
Set-TextCell 'E19' '  -1.05%  '

# Row 20: 'BitcoinCash'
Set-TextCell 'D20' '445.20'
Set-TextCell 'E20' '  -0.86%  '

# Row 21: 'Chainlink'
Set-TextCell 'D21' '13.73'
Set-TextCell 'E21' '  +0.48%  '

# Row 22: 'Polygon'
Set-TextCell 'D22' '0.674'
Set-TextCell 'E22' '  -0.86%  '

# Row 23: 'Uniswap'
Set-TextCell 'D23' '7.18'
Set-TextCell 'E23' '  -1.77%  '

# Row 24: 'Litecoin'
Set-TextCell 'D24' '82.23'
Set-TextCell 'E24' '  +1.34%  '

# Row 25: 'Fetch.AI'
Set-TextCell 'E25' '  -0.47%  '

# Row 26: 'InternetComputer(DFINITY)'
Set-TextCell 'D26' '12.17'
Set-TextCell 'E26' '  -0.66%  '

# Row 27: 'Dai' -> 'RenderToken'
Set-TextCell 'B27' 'RenderToken'
Set-TextCell 'C27' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D27' '10.04'
Set-TextCell 'E27' '  -8.93%  '

# Row 28: 'RenderToken' -> 'Dai'
Set-TextCell 'B28' 'Dai'
Set-TextCell 'C28' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 'D28' '1.00'
Set-TextCell 'E28' '  +0.02%  '

# Row 29: 'NEARProtocol'
Set-TextCell 'D29' '8.09'
Set-TextCell 'E29' '  +3.63%  '

# Row 30: 'ImmutableX'
Set-TextCell 'D30' '2.33'
Set-TextCell 'E30' '  -1.56%  '

# Row 31: 'PancakeSwap'
Set-TextCell 'D31' '2.57'
Set-TextCell 'E31' '  -0.11%  '

# Row 32: 'PEPE'
Set-TextCell 'D32' '0.0₃0974'
Set-TextCell 'E32' '  -10.18%  '

# Row 33: 'EthereumClassic'
Set-TextCell 'D33' '27.16'
Set-TextCell 'E33' '  +2.05%  '

# Row 34: 'Hedera'
Set-TextCell 'D34' '0.110'
Set-TextCell 'E34' '  -0.12%  '

# Row 35: 'FirstDigitalUSD'
Set-TextCell 'D35' '0.998'
Set-TextCell 'E35' '  -0.07%  '

# Row 36: 'Mantle'
Set-TextCell 'D36' '0.975'
Set-TextCell 'E36' '  -1.00%  '

# Row 37: 'Filecoin'
Set-TextCell 'D37' '5.69'
Set-TextCell 'E37' '  +0.35%  '

# Row 38: 'OKB'
Set-TextCell 'D38' '49.17'
Set-TextCell 'E38' '  +0.26%  '

# Row 39: 'Stacks'
Set-TextCell 'E39' '  -6.84%  '

# Row 40: 'TheGraph'
Set-TextCell 'E40' '  +1.56%  '

# Row 41: 'Arweave' -> 'Kaspa'
Set-TextCell 'B41' 'Kaspa'
Set-TextCell 'C41' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell 'D41' '0.119'
Set-TextCell 'E41' '  -1.51%  '

# Row 42: 'dogwifhat' -> 'Arweave'
Set-TextCell 'B42' 'Arweave'
Set-TextCell 'C42' 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextCell 'D42' '43.00'
Set-TextCell 'E42' '  -2.99%  '

# Row 43: 'Kaspa' -> 'dogwifhat'
Set-TextCell 'B43' 'dogwifhat'
Set-TextCell 'C43' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell 'D43' '2.80'
Set-TextCell 'E43' '  -7.07%  '

# Row 44: 'Cosmos'
Set-TextCell 'D44' '8.37'
Set-TextCell 'E44' '  -0.39%  '

# Row 45: 'Bittensor'
Set-TextCell 'D45' '382.91'
Set-TextCell 'E45' '  -1.36%  '

# Row 46: 'VeChain'
Set-TextCell 'D46' '0.0353'
Set-TextCell 'E46' '  +1.01%  '

# Row 47: 'Maker'
Set-TextCell 'D47' '2.715.00'
Set-TextCell 'E47' '  -2.65%  '

# Row 48: 'Monero'
Set-TextCell 'D48' '130.32'
Set-TextCell 'E48' '  -3.60%  '

# Row 49: 'USDe'
Set-TextCell 'E49' '  +0.04%  '

# Row 50: 'Stellar'
Set-TextCell 'E50' '  +0.73%  '

# Row 51: 'InjectiveProtocol'
Set-TextCell 'D51' '23.11'
Set-TextCell 'E51' '  +0.20%  '
